# Fill in the "Completion" column (D) of the self-assessment rubric table
# with the actual completion scores, and leave the selection where the
# user finished typing (cell D22, just past the last filled row D21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = 0.75
# D6 intentionally left blank (no change)
$ws.Range("D7").Value = 0.25
$ws.Range("D8").Value = 0.25

$ws.Range("D9").Value = 0.5
$ws.Range("D9").Font.Bold = $true
# D10 intentionally left blank (no change)

$ws.Range("D11").Value = 0
$ws.Range("D12").Value = 0.5
$ws.Range("D13").Value = 0.25
$ws.Range("D14").Value = 0.5
$ws.Range("D15").Value = 1.25
$ws.Range("D16").Value = 0.25
$ws.Range("D17").Value = 0.5
$ws.Range("D18").Value = 0.5
$ws.Range("D19").Value = 0.5
$ws.Range("D20").Value = 0

# Row 21 is the bold "Total" row
$ws.Range("D21").Value = 6
$ws.Range("D21").Font.Bold = $true

$ws.Range("D22").Select()
